$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("TestSheet_1")

# Remember original MW column values (I2:I13) before we clear them
$mwValues = $src.Range("I2:I13").Value

# Duplicate TestSheet_1 and place the copy right after it; rename to TestSheet_1_withMW
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "TestSheet_1_withMW"

# Clear the Molecular Weight values from the original TestSheet_1 (keep the formatting)
$src.Range("I2:I13").ClearContents()
